# Update the cryptos list with the latest price and volume(1h) figures.
# D-column prices are plain text (locale-formatted, e.g. thousand separators
# using '.' and fixed trailing zeros) so we force text storage before/while
# assigning to avoid Excel auto-converting them to numbers, then reset the
# style back to Normal so no extra formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.375.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.301.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.652"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.860.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("E13").Value = "  -5.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.357.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.282.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("E23").Value = "  +0.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.425.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.513"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.198"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.38%  "

$ws.Range("E28").Value = "  -3.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.804.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.776"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0664"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "321.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0270"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("E51").Value = "  +6.29%  "
